$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Enter the new test-case rows in the order the original author likely
# typed them (A,C,B for row 42; C,B,A for row 43) so the shared-string table
# is rebuilt in the same order as the target workbook. ---
$ws.Range("A42").Value = "Profile41"
$ws.Range("C42").Value = "Verify that HCR having Badge append with his name in Profile page"
$ws.Range("B42").Value = "OPQA-2726"
$ws.Range("D42").Value = "Y"
$ws.Range("C43").Value = "Verify that user is able to edit first name and last name from his own profile page."
$ws.Range("B43").Value = "OPQA-2679"
$ws.Range("A43").Value = "Profile42"
$ws.Range("D43").Value = "Y"

# --- Apply matching cell styles by copying the formatting of the existing
# analogous cells (keeps reusing the workbook's existing style table instead
# of minting new cellXfs entries). ---
$ws.Range("A41").Copy()
$ws.Range("A42").PasteSpecial(-4122)

$ws.Range("B41").Copy()
$ws.Range("B42").PasteSpecial(-4122)
$ws.Range("B43").PasteSpecial(-4122)

$ws.Range("C41").Copy()
$ws.Range("C42").PasteSpecial(-4122)
$ws.Range("C43").PasteSpecial(-4122)

$ws.Range("D41").Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("D43").PasteSpecial(-4122)

$ws.Range("E41").Copy()
$ws.Range("E42").PasteSpecial(-4122)
$ws.Range("E43").PasteSpecial(-4122)

# Row 43's TCID cell (A43) carries the plain bordered style used by the
# C column rather than the usual A-column style.
$ws.Range("C21").Copy()
$ws.Range("A43").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Hyperlinks to the Jira issues referenced by the two new test cases ---
$ws.Hyperlinks.Add($ws.Range("B42"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-2726", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-2726")
$ws.Range("B42").Value = "OPQA-2726"

$ws.Hyperlinks.Add($ws.Range("B43"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-2679", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-2679")
$ws.Range("B43").Value = "OPQA-2679"

# Hyperlinks.Add auto-applies Excel's built-in "Hyperlink" look to the cell;
# restore the original (non-hyperlink-blue) formatting used by the rest of
# the sheet, as the test-case tracker doesn't style these cells specially.
$ws.Range("B41").Copy()
$ws.Range("B42").PasteSpecial(-4122)
$ws.Range("B43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Drop the now-unused built-in Hyperlink cell style that Excel auto-registers
# on the workbook the first time a hyperlink is added.
try {
  $wb.Styles.Item("Hyperlink").Delete()
} catch {
}

# --- Sheet view: scrolled down with C45 selected (matches the new bottom
# of the test-case list) ---
$ws.Range("C45").Select()

Write-Host "done"
